# Repull data: update column F (dSF) values for specific rows on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = 0
    9  = -2
    11 = -7
    13 = -12
    14 = -7
    15 = 0
    20 = -10
    24 = 6
    26 = -5
    29 = -7
    34 = -2
    36 = -2
    39 = 1
    42 = 5
    43 = -3
    44 = -2
    45 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
